$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column H (SAT) for rows 2 through 20 from 324 to 272
$ws.Range("H2:H20").Value = 272
